# kim_product_4_data.xlsx — dataloader modified and environment data only
# train/infer supported.
#
# Column A held Excel date serials (e.g. 43227 == 2018-05-07) displayed
# via a custom "YYYY-MM-DD HH:MM:SS" number format. The new dataloader
# expects plain numeric dates in YYYYMMDD form (e.g. 20180507) with the
# default/general number format instead, so convert every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)

    # Value2 returns the raw underlying Excel serial number; Value would
    # come back as a formatted .NET DateTime because of the cell's
    # current date number format.
    $serial = $cell.Value2()

    if ($serial -eq $null) {
        continue
    }

    # Excel serial date -> proleptic date, using the 1899-12-30 epoch.
    $epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
    $date = $epoch.AddDays([double]$serial)
    $ymd = [int]$date.ToString("yyyyMMdd")

    $cell.Value = $ymd

    # Drop the custom date/time number format entirely and fall back to
    # the sheet's default (unstyled) look, matching the new plain-integer
    # data convention.
    $cell.Style = "Normal"
}
